$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 60428
$ws.Range("E2").Value = 1566
$ws.Range("F2").Value = 1566
$ws.Range("G2").Value = 1546
$ws.Range("H2").Value = 1149
$ws.Range("I2").Value = 1145
$ws.Range("J2").Value = 3
$ws.Range("K2").Value = 130474
$ws.Range("L2").Value = 116452
$ws.Range("M2").Value = 14022
$ws.Range("N2").Value = 13976
$ws.Range("O2").Value = 46
$ws.Range("P2").Value = 530
$ws.Range("Q2").Value = 817
$ws.Range("R2").Value = -368
$ws.Range("S2").Value = 320
$ws.Range("T2").Value = 178
$ws.Range("V2").Value = 2452
$ws.Range("W2").Value = 2.59
$ws.Range("X2").Value = 1.9
$ws.Range("Y2").Value = 9.43
$ws.Range("Z2").Value = 0.95
$ws.Range("AA2").Value = 830.51
$ws.Range("AB2").Value = 2569.21
$ws.Range("AC2").Value = 1119
$ws.Range("AD2").Value = 11.13
$ws.Range("AE2").Value = 13291
$ws.Range("AF2").Value = 0.94
$ws.Range("AG2").Value = 380
$ws.Range("AH2").Value = 3.05
$ws.Range("AI2").Value = 34.88
$ws.Range("AJ2").Value = 105963000

# Row 3
$ws.Range("D3").Value = 67930
$ws.Range("E3").Value = 2247
$ws.Range("F3").Value = 2247
$ws.Range("G3").Value = 2222
$ws.Range("H3").Value = 1690
$ws.Range("I3").Value = 1685
$ws.Range("J3").Value = 5
$ws.Range("K3").Value = 148325
$ws.Range("L3").Value = 133513
$ws.Range("M3").Value = 14811
$ws.Range("N3").Value = 14762
$ws.Range("O3").Value = 49
$ws.Range("P3").Value = 530
$ws.Range("Q3").Value = 9596
$ws.Range("R3").Value = -9357
$ws.Range("S3").Value = -53
$ws.Range("T3").Value = 99
$ws.Range("V3").Value = 2852
$ws.Range("W3").Value = 3.31
$ws.Range("X3").Value = 2.49
$ws.Range("Y3").Value = 11.76
$ws.Range("Z3").Value = 1.21
$ws.Range("AA3").Value = 901.43
$ws.Range("AB3").Value = 2706.98
$ws.Range("AC3").Value = 1590
$ws.Range("AD3").Value = 10.16
$ws.Range("AE3").Value = 13985
$ws.Range("AF3").Value = 1.15
$ws.Range("AG3").Value = 570
$ws.Range("AH3").Value = 3.53
$ws.Range("AI3").Value = 35.7
$ws.Range("AJ3").Value = 105963000

# Row 4
$ws.Range("D4").Value = 71520
$ws.Range("E4").Value = 3143
$ws.Range("F4").Value = 3143
$ws.Range("G4").Value = 3105
$ws.Range("H4").Value = 2372
$ws.Range("I4").Value = 2365
$ws.Range("J4").Value = 7
$ws.Range("K4").Value = 165737
$ws.Range("L4").Value = 149351
$ws.Range("M4").Value = 16385
$ws.Range("N4").Value = 16326
$ws.Range("O4").Value = 60
$ws.Range("P4").Value = 552
$ws.Range("Q4").Value = 8700
$ws.Range("R4").Value = -8430
$ws.Range("S4").Value = 0
$ws.Range("T4").Value = 62
$ws.Range("V4").Value = 2853
$ws.Range("W4").Value = 4.4
$ws.Range("X4").Value = 3.32
$ws.Range("Y4").Value = 15.26
$ws.Range("Z4").Value = 1.51
$ws.Range("AA4").Value = 911.49
$ws.Range("AB4").Value = 2889.76
$ws.Range("AC4").Value = 2175
$ws.Range("AD4").Value = 7.04
$ws.Range("AE4").Value = 14893
$ws.Range("AF4").Value = 1.03
$ws.Range("AG4").Value = 830
$ws.Range("AH4").Value = 5.42
$ws.Range("AI4").Value = 38.47
$ws.Range("AJ4").Value = 110338000

# Row 5
$ws.Range("D5").Value = 79335
$ws.Range("E5").Value = 5136
$ws.Range("F5").Value = 5136
$ws.Range("G5").Value = 5125
$ws.Range("H5").Value = 3846
$ws.Range("I5").Value = 3838
$ws.Range("J5").Value = 8
$ws.Range("K5").Value = 181524
$ws.Range("L5").Value = 163732
$ws.Range("M5").Value = 17793
$ws.Range("N5").Value = 17732
$ws.Range("O5").Value = 60
$ws.Range("P5").Value = 552
$ws.Range("Q5").Value = 22754
$ws.Range("R5").Value = -20061
$ws.Range("S5").Value = -960
$ws.Range("T5").Value = 91
$ws.Range("V5").Value = 2855
$ws.Range("W5").Value = 6.47
$ws.Range("X5").Value = 4.85
$ws.Range("Y5").Value = 22.59
$ws.Range("Z5").Value = 2.21
$ws.Range("AA5").Value = 920.23
$ws.Range("AB5").Value = 3157.77
$ws.Range("AC5").Value = 3479
$ws.Range("AD5").Value = 6.76
$ws.Range("AE5").Value = 16237
$ws.Range("AF5").Value = 1.45
$ws.Range("AG5").Value = 1140
$ws.Range("AH5").Value = 4.85
$ws.Range("AI5").Value = 32.44
$ws.Range("AJ5").Value = 110338000

# Row 6
$ws.Range("D6").Value = 84182
$ws.Range("E6").Value = 3127
$ws.Range("F6").Value = 3127
$ws.Range("G6").Value = 3148
$ws.Range("H6").Value = 2347
$ws.Range("I6").Value = 2338
$ws.Range("K6").Value = 204788
$ws.Range("L6").Value = 181956
$ws.Range("M6").Value = 22832
$ws.Range("N6").Value = 22765
$ws.Range("P6").Value = 568
$ws.Range("Q6").Value = 13225
$ws.Range("R6").Value = -13621
$ws.Range("S6").Value = 297
$ws.Range("T6").Value = 42
$ws.Range("V6").Value = 3853
$ws.Range("W6").Value = 3.72
$ws.Range("X6").Value = 2.79
$ws.Range("Y6").Value = 11.59
$ws.Range("Z6").Value = 1.21
$ws.Range("AA6").Value = 796.93
$ws.Range("AB6").Value = 3976.72
$ws.Range("AC6").Value = 2088
$ws.Range("AD6").Value = 10.46
$ws.Range("AE6").Value = 20348
$ws.Range("AF6").Value = 1.07
$ws.Range("AG6").Value = 820
$ws.Range("AH6").Value = 3.75
$ws.Range("AI6").Value = 39.24
$ws.Range("AJ6").Value = 113680000

# Row 7
$ws.Range("D7").Value = 76781
$ws.Range("E7").Value = 4087
$ws.Range("G7").Value = 4164
$ws.Range("H7").Value = 2821
$ws.Range("I7").Value = 2770
$ws.Range("K7").Value = 230842
$ws.Range("L7").Value = 202421
$ws.Range("M7").Value = 28335
$ws.Range("N7").Value = 29190
$ws.Range("P7").Value = 568
$ws.Range("W7").Value = 5.32
$ws.Range("X7").Value = 3.67
$ws.Range("Y7").Value = 10.66
$ws.Range("Z7").Value = 1.29
$ws.Range("AA7").Value = 714.38
$ws.Range("AC7").Value = 2437
$ws.Range("AD7").Value = 6.5
$ws.Range("AE7").Value = 26183
$ws.Range("AF7").Value = 0.61
$ws.Range("AG7").Value = 821
$ws.Range("AH7").Value = 5.18
$ws.Range("AI7").Value = 33.68

# Row 8
$ws.Range("D8").Value = 87781
$ws.Range("E8").Value = 3999
$ws.Range("G8").Value = 3217
$ws.Range("H8").Value = 2321
$ws.Range("I8").Value = 2195
$ws.Range("K8").Value = 258596
$ws.Range("L8").Value = 227411
$ws.Range("M8").Value = 29712
$ws.Range("N8").Value = 30510
$ws.Range("P8").Value = 568
$ws.Range("W8").Value = 4.56
$ws.Range("X8").Value = 2.64
$ws.Range("Y8").Value = 7.35
$ws.Range("Z8").Value = 0.95
$ws.Range("AA8").Value = 765.38
$ws.Range("AC8").Value = 1931
$ws.Range("AD8").Value = 8.13
$ws.Range("AE8").Value = 27367
$ws.Range("AF8").Value = 0.57
$ws.Range("AG8").Value = 792
$ws.Range("AH8").Value = 5.04
$ws.Range("AI8").Value = 40.99

# Row 9
$ws.Range("D9").Value = 101810
$ws.Range("E9").Value = 5253
$ws.Range("G9").Value = 3905
$ws.Range("H9").Value = 2822
$ws.Range("I9").Value = 2365
$ws.Range("K9").Value = 287797
$ws.Range("L9").Value = 256178
$ws.Range("M9").Value = 31618
$ws.Range("N9").Value = 31995
$ws.Range("P9").Value = 568
$ws.Range("W9").Value = 5.16
$ws.Range("X9").Value = 2.77
$ws.Range("Y9").Value = 7.57
$ws.Range("Z9").Value = 1.03
$ws.Range("AA9").Value = 810.22
$ws.Range("AC9").Value = 2080
$ws.Range("AD9").Value = 7.55
$ws.Range("AE9").Value = 28699
$ws.Range("AF9").Value = 0.55
$ws.Range("AG9").Value = 876
$ws.Range("AH9").Value = 22.63
$ws.Range("AI9").Value = 42.11

# Clear removed cells
$ws.Range("U2").ClearContents()
$ws.Range("U3").ClearContents()
$ws.Range("U4").ClearContents()
$ws.Range("U5").ClearContents()
$ws.Range("U6").ClearContents()
$ws.Range("Q7").ClearContents()
$ws.Range("R7").ClearContents()
$ws.Range("S7").ClearContents()
$ws.Range("T7").ClearContents()
$ws.Range("U7").ClearContents()
$ws.Range("Q8").ClearContents()
$ws.Range("R8").ClearContents()
$ws.Range("S8").ClearContents()
$ws.Range("T8").ClearContents()
$ws.Range("U8").ClearContents()
$ws.Range("Q9").ClearContents()
$ws.Range("R9").ClearContents()
$ws.Range("S9").ClearContents()
$ws.Range("T9").ClearContents()
$ws.Range("U9").ClearContents()